# Fix Training Data Issue
# The "Date" column (BF) on this sheet was stamped with the wrong value
# ("4-22-2012-13", a mangled concatenation of month-day and season) for
# every data row. It should read the actual game date "2013-04-22".
#
# Excel normally auto-parses a plain "YYYY-MM-DD" string assignment into a
# date serial number, which would silently change both the stored value and
# the cell's displayed text (e.g. "4/22/13"). To keep the value as the exact
# literal text "2013-04-22" (matching the original cell's plain-text style),
# we build the text in an unused scratch cell via a formula (which returns a
# text result, bypassing Excel's literal-entry date autodetection) and copy
# that already-evaluated text into each target cell with a values-only paste.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$correctDate = "2013-04-22"

# Scratch cell well outside the used range (A1:BF31) so it never collides
# with real data and can be cleared afterwards without leaving a trace.
$helper = $ws.Range("ZZ1")
$helper.Formula = '="' + $correctDate + '"'
$helper.Copy()

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    if ($cell.Text -eq "4-22-2012-13") {
        $cell.PasteSpecial(-4163)  # xlPasteValues
    }
}

# Clean up the scratch cell and clipboard marquee.
$helper.ClearContents()
$excel.CutCopyMode = $false
